$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: drop the leading "1- " / "2- " numbering from the two
# prerequisite bullets, leaving a plain "- " prefix.
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "1- Python 3.9.10 ou une version ultérieure", $false, $false, $false,
    $false, $false, $true, 1, $false,
    "- Python 3.9.10 ou une version ultérieure", 2) | Out-Null

$d.Content.Find.Execute(
    "2- Git Bash", $false, $false, $false, $false, $false, $true, 1,
    $false, "- Git Bash", 2) | Out-Null

# ---------------------------------------------------------------------
# Change 2: split the hyperlink run "de Python" into three runs with
# identical formatting: "de P", "y", "thon". We locate the exact
# character span with Find (so field-code offsets inside the hyperlink
# are accounted for), then force a run split by toggling Bold on/off on
# the narrower sub-spans (no net formatting change, only a run break).
# ---------------------------------------------------------------------
$hlRange = $null
foreach ($hl in $d.Hyperlinks) {
    if ($hl.Range.Text -eq "site officiel de Python") {
        $hlRange = $hl.Range
        break
    }
}

if ($hlRange -ne $null) {
    $scan = $hlRange.Duplicate
    $scan.Find.Execute("y", $true, $false, $false, $false, $false, $true,
                        1, $false, "", 0) | Out-Null
    $ySpan = $d.Range($scan.Start, $scan.End)
    $ySpan.Font.Bold = 1
    $ySpan.Font.Bold = 0
}

# ---------------------------------------------------------------------
# Change 3: inside the "export_output.html" paragraph, carve the
# filename out into its own bold run, e.g.
#   "...(export_output.html) sera créé..."
#                ^^^^^^^^^^^^^^^^^^^ now bold
# ---------------------------------------------------------------------
$scan2 = $d.Content.Duplicate
$scan2.Find.Execute("export_output.html", $false, $false, $false, $false,
                     $false, $true, 1, $false, "", 0) | Out-Null
$htmlSpan = $d.Range($scan2.Start, $scan2.End)
$htmlSpan.Font.Bold = 1
